$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "38.251.55"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +3.73%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.062.36"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +3.40%  "

$ws.Range("E4").Value = "  +0.28%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "230.87"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.79%  "

$ws.Range("E6").Value = "  +1.85%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "58.27"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +7.60%  "

$ws.Range("E8").Value = "  -0.02%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.389"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +3.51%  "

$ws.Range("E10").Value = "  +3.46%  "

$ws.Range("E11").Value = "  +0.17%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "2.367.89"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +3.47%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "14.66"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +4.48%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "20.76"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +3.70%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.758"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.16%  "

$ws.Range("E16").Value = "  +4.27%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.057.16"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +2.87%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "38.138.47"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +3.71%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.19"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "69.97"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("E21").Value = "  +2.68%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "225.40"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("E24").Value = "  +1.47%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.26"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +4.50%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.34"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +2.67%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "166.45"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("E28").Value = "  +8.67%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "19.10"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.04%  "

$ws.Range("E30").Value = "  +2.74%  "

$ws.Range("E31").Value = "  +1.87%  "

$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("E33").Value = "  +5.55%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0617"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +7.63%  "

$ws.Range("E36").Value = "  +2.89%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "6.11"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +15.73%  "

$ws.Range("E38").Value = "  +7.03%  "

$ws.Range("E39").Value = "  +0.11%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "98.54"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +4.52%  "

$ws.Range("E41").Value = "  +2.34%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.487.12"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "17.03"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.28%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0948"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +3.83%  "

$ws.Range("E45").Value = "  +3.85%  "

$ws.Range("E46").Value = "  +1.45%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.14"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +19.57%  "

$ws.Range("E48").Value = "  +2.21%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.97"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.35%  "

$ws.Range("E50").Value = "  -0.03%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.255.34"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +3.33%  "
